# Update the cryptos price/volume table with refreshed values.
# For D-column cells whose new value looks like a plain number (e.g. "7.15"),
# force the cell to Text format first so Excel keeps storing it as text
# (matching the workbook's existing inlineStr string cells) instead of
# silently converting it to a floating point number and losing formatting
# such as trailing zeros (e.g. "36.70" -> 36.7).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = '58.230.36'
$ws.Range("E2").Value = '  -0.41%  '
$ws.Range("D3").Value = '3.140.05'
$ws.Range("E3").Value = '  +1.72%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '530.54'
$ws.Range("E5").Value = '  +1.35%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.45'
$ws.Range("E6").Value = '  -0.58%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = '3.137.54'
$ws.Range("E8").Value = '  +1.66%  '
$ws.Range("E9").Value = '  +1.19%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.15'
$ws.Range("E10").Value = '  -2.48%  '
$ws.Range("E11").Value = '  +0.51%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.393'
$ws.Range("E12").Value = '  +2.27%  '
$ws.Range("D13").Value = '3.683.90'
$ws.Range("E13").Value = '  +1.78%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.135'
$ws.Range("E14").Value = '  +3.45%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '25.61'
$ws.Range("E15").Value = '  -4.31%  '
$ws.Range("E16").Value = '  -0.32%  '
$ws.Range("D17").Value = '58.295.65'
$ws.Range("E17").Value = '  -0.44%  '
$ws.Range("D18").Value = '3.133.41'
$ws.Range("E18").Value = '  +1.61%  '
$ws.Range("E19").Value = '  -0.47%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.82'
$ws.Range("E20").Value = '  -0.85%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.98'
$ws.Range("E21").Value = '  -1.35%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '343.31'
$ws.Range("E22").Value = '  +0.35%  '
$ws.Range("E23").Value = '  -0.22%  '
$ws.Range("E24").Value = '  +1.81%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '67.64'
$ws.Range("E25").Value = '  +2.83%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.170'
$ws.Range("E26").Value = '  -0.48%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  -0.11%  '
$ws.Range("D28").Value = '0.0₃0932'
$ws.Range("E28").Value = '  +2.00%  '
$ws.Range("E29").Value = '  +0.02%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.37'
$ws.Range("E30").Value = '  +2.35%  '
$ws.Range("E31").Value = '  -2.38%  '
$ws.Range("E32").Value = '  +1.76%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '21.15'
$ws.Range("E33").Value = '  +0.82%  '
$ws.Range("E34").Value = '  -0.39%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '158.70'
$ws.Range("E35").Value = '  +2.80%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.78'
$ws.Range("E36").Value = '  +4.12%  '
$ws.Range("E37").Value = '  +2.52%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '26.33'
$ws.Range("E38").Value = '  -2.19%  '
$ws.Range("E39").Value = '  -4.47%  '
$ws.Range("E40").Value = '  +12.32%  '
$ws.Range("E41").Value = '  -1.42%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.707'
$ws.Range("E42").Value = '  +5.73%  '
$ws.Range("E43").Value = '  +2.58%  '
$ws.Range("D44").Value = '3.184.32'
$ws.Range("E44").Value = '  +1.64%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '36.70'
$ws.Range("E46").Value = '  -0.02%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0265'
$ws.Range("E47").Value = '  +3.42%  '
$ws.Range("D48").Value = '2.274.50'
$ws.Range("E48").Value = '  +0.02%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.01'
$ws.Range("E49").Value = '  +5.28%  '
$ws.Range("E50").Value = '  +1.94%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '20.66'
$ws.Range("E51").Value = '  -0.45%  '
